$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 176 (shifts existing rows 176..296 down to 177..297)
$ws.Rows(176).Insert()

# Populate the newly inserted row with its data
$ws.Range("A176").Value = 11
$ws.Range("B176").Value = "Vega Monumental Concepción"
$ws.Range("C176").Value = "Bíobío"
$ws.Range("D176").Value = 44673
$ws.Range("E176").Value = 8
$ws.Range("F176").Value = 100114014
$ws.Range("G176").Value = "Betarraga"
$ws.Range("H176").Value = "Sin especificar"
$ws.Range("I176").Value = "Primera"
$ws.Range("J176").Value = 550
$ws.Range("K176").Value = 600
$ws.Range("L176").Value = 650
$ws.Range("M176").Value = 627
$ws.Range("N176").Value = "`$/paquete 5 unidades"
$ws.Range("O176").Value = "Región Metropolitana"
$ws.Range("P176").Value = 125
$ws.Range("Q176").Value = 5
$ws.Range("R176").Value = "Hortaliza"
